$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.011503338813782
$ws.Range("B1").Value = 1.42776894569397
$ws.Range("D1").Value = 1.729949235916138
$ws.Range("E1").Value = 1.035499930381775
